# 2573-...-ADD-VAR-INST-CASH: "Loan RBI, Variable Instalments"
# A new (blank) column is inserted into the "Repayment schedule" sheet, just
# before the existing "Late" column, to make room for a Variable Instalments
# related value; the trailing columns (Late / heading / Outstanding) shift
# one place to the right. The "Repayment schedule" tab is also left as the
# active tab/selection (S9) when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (pushes old N:P -> O:Q).
[void]$ws.Columns("N:N").Insert()

# The inserted column keeps the same width as column M (11 "characters"),
# but without the bestFit flag since it now starts out blank.
$ws.Columns("N:N").ColumnWidth = 10.166666666666666

# Make "Repayment schedule" the active/selected sheet, with S9 selected -
# matches the workbook's saved view state after the edit.
[void]$ws.Activate()
[void]$ws.Range("S9").Select()
